$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ------------------------------------------------------------------
# 1) Swap the header fills of columns B and C (B1 <-> C1), using a
#    scratch cell so neither format is lost while swapping.
# ------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)     # E1 scratch now has C1's (khaki) format

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial($xlPasteFormats)      # C1 now has B1's (indianred) format

$ws.Range("E1").Copy()
$ws.Range("B1").PasteSpecial($xlPasteFormats)      # B1 now has the original C1 (khaki) format

$ws.Range("E1").Clear()

# ------------------------------------------------------------------
# 2) Move the highlighted (salmon) fill from B2:B3 to C2:C3, and
#    clear the highlight from B2:B3 back to the default style.
# ------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("C2:C3").PasteSpecial($xlPasteFormats)

$ws.Range("B2:B3").ClearFormats()

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Re-order the columns (B <-> C) and the data rows (2 <-> 3),
#    and apply the new header labels. Values are written last so
#    the copy/paste steps above are unaffected by the new content.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Coluna C"
$ws.Range("C1").Value = "Coluna B"

$ws.Range("A2").Value = "Linha A3"
$ws.Range("B2").Value = "Linha C3"
$ws.Range("C2").ClearContents()

$ws.Range("A3").Value = "Linha A2"
$ws.Range("B3").Value = "Linha C2"
$ws.Range("C3").Value = "Linha BB"

$ws.Range("B4").Value = "Linha C4"
$ws.Range("C4").Value = "Linha B4"
